$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text (conversion rates) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 3.36 = 12818.79 pesos", "1000 Bs = 3.34 = 12757.53 pesos")
$text = $text.Replace("12818.79 pesos = 3.36 = 976.58 Bs", "12757.53 pesos = 3.33 = 969.36 Bs")
$cell.Value = $text

# --- Update tasas sheet numeric values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 299
$wsTasas.Range("O10").Value = 3814.5
$wsTasas.Range("N12").Value = 3825.83
$wsTasas.Range("O12").Value = 290.7
